$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 153, pushing the existing rows 153-278 down
# to 154-279 (matches the dimension growing from A1:R278 to A1:R279).
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row 153 with the new weekly price record.
$ws.Range("A153").Value = 7
$ws.Range("B153").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C153").Value = "Ñuble"
$ws.Range("D153").Value = 44942
$ws.Range("E153").Value = 16
$ws.Range("F153").Value = 100112017
$ws.Range("G153").Value = "Apio"
$ws.Range("H153").Value = "Americana (o)"
$ws.Range("I153").Value = "Primera"
$ws.Range("J153").Value = 60
$ws.Range("K153").Value = 10000
$ws.Range("L153").Value = 10000
$ws.Range("M153").Value = 10000
$ws.Range("N153").Value = "$/docena de matas"
$ws.Range("O153").Value = "Provincia del Elquí"
$ws.Range("P153").Value = 1667
$ws.Range("Q153").Value = 6
$ws.Range("R153").Value = "Hortaliza"
